$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 534.6510700667664
$ws.Range("C2").Value = 377.7946475571405
$ws.Range("D2").Value = 324.051286191299
$ws.Range("E2").Value = 294.1235762200708
$ws.Range("B3").Value = 629.9092690074646
$ws.Range("C3").Value = 447.0569051666072
$ws.Range("D3").Value = 381.4389287730951
$ws.Range("E3").Value = 348.9452654021009
$ws.Range("B4").Value = 598.1958735930043
$ws.Range("C4").Value = 424.1184117712934
$ws.Range("D4").Value = 362.0099379655671
$ws.Range("E4").Value = 331.7329920809061
$ws.Range("B5").Value = 412.9368477255957
$ws.Range("C5").Value = 293.7040883048616
$ws.Range("D5").Value = 247.5866598743755
$ws.Range("E5").Value = 229.1525758334074
$ws.Range("B6").Value = 368.1656625667244
$ws.Range("C6").Value = 260.2542893699746
$ws.Range("D6").Value = 220.7089877195534
$ws.Range("E6").Value = 201.5736754565277
$ws.Range("B7").Value = 38.240897399423
$ws.Range("C7").Value = 27.00912438045556
$ws.Range("D7").Value = 23.14028073250289
$ws.Range("E7").Value = 21.07467748317314
$ws.Range("B8").Value = 2108.712530461781
$ws.Range("C8").Value = 1492.364264413481
$ws.Range("D8").Value = 1288.178366683002
$ws.Range("E8").Value = 1159.863934197644
$ws.Range("B9").Value = 529.2628919600966
$ws.Range("C9").Value = 376.1099035608317
$ws.Range("D9").Value = 320.922848532379
$ws.Range("E9").Value = 293.7669522327109
$ws.Range("B10").Value = 196.6119278183982
$ws.Range("C10").Value = 136.2395603756289
$ws.Range("D10").Value = 119.6487339999811
$ws.Range("E10").Value = 109.4132806796197
$ws.Range("B11").Value = 35.10725127135109
$ws.Range("C11").Value = 22.9913519096603
$ws.Range("D11").Value = 19.94071111510026
$ws.Range("E11").Value = 19.71785220273538
$ws.Range("B12").Value = 72.76555853056067
$ws.Range("C12").Value = 51.23873475999182
$ws.Range("D12").Value = 43.30216753786912
$ws.Range("E12").Value = 38.97354084412795
$ws.Range("B13").Value = 109.6583159454891
$ws.Range("C13").Value = 75.12354742175511
$ws.Range("D13").Value = 66.16126303641018
$ws.Range("E13").Value = 61.19905936475248
